# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.061.02"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.664.16"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'310.34"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Value = "'51.66"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").Value = "'1.374"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.08518"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "'24.19"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "'7.246"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "'8.018"
$ws.Range("E15").Value = "  +6.94%  "
$ws.Range("D16").Value = "'0.00001321"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "1.661.73"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "'94.92"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "'0.06996"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'20.04"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'7.013"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'13.76"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "24.062.16"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'3.164"
$ws.Range("E25").Value = "  +9.81%  "
$ws.Range("D26").Value = "'2.498"
$ws.Range("E26").Value = "  +3.67%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "'153.89"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").Value = "'141.35"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'5.311"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "'7.849"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'2.482"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "1.845.48"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'1.055"
$ws.Range("E34").Value = "  +8.71%  "
$ws.Range("D35").Value = "'0.08194"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("D36").Value = "'0.03026"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("D37").Value = "'11.29"
$ws.Range("E37").Value = "  +9.19%  "
$ws.Range("D38").Value = "'6.740"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "'0.2724"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").Value = "'0.09177"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'13.76"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("D42").Value = "'0.7627"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'1.433"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'16.60"
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").Value = "'0.7053"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").Value = "'2.516"
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("D47").Value = "'4.104"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.08344"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'135.60"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'1.246"
$ws.Range("E51").Value = "  -1.54%  "
